# Weekly data refresh: a new price record (week of 2021-09-16) is added
# for "Vega Monumental Concepción - Papa" / Asterix / 1a (guarda).
# It belongs at row 81 (sheet is sorted by date descending), so insert a
# new row there; every subsequent record shifts down by one row
# (old row 81 -> 82, ..., old row 111 -> 112).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(81).Insert()

$ws.Cells.Item(81, 1).Value  = 11
$ws.Cells.Item(81, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(81, 3).Value  = "Bíobío"
$ws.Cells.Item(81, 4).Value  = 44455
$ws.Cells.Item(81, 5).Value  = 8
$ws.Cells.Item(81, 6).Value  = 100114001
$ws.Cells.Item(81, 7).Value  = "Papa"
$ws.Cells.Item(81, 8).Value  = "Asterix"
$ws.Cells.Item(81, 9).Value  = "1a (guarda)"
$ws.Cells.Item(81, 10).Value = 2000
$ws.Cells.Item(81, 11).Value = 8500
$ws.Cells.Item(81, 12).Value = 9000
$ws.Cells.Item(81, 13).Value = 8750
$ws.Cells.Item(81, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(81, 15).Value = "Provincia de Arauco"
$ws.Cells.Item(81, 16).Value = 350
$ws.Cells.Item(81, 17).Value = 25
$ws.Cells.Item(81, 18).Value = "Hortaliza"
